$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap species data between row 4 and row 7 (columns A,B,D,E,F,G,H) ---
$row4A = $ws.Range("A4").Value2
$row4B = $ws.Range("B4").Value2
$row4D = $ws.Range("D4").Value2
$row4E = $ws.Range("E4").Value2
$row4F = $ws.Range("F4").Value2
$row4G = $ws.Range("G4").Value2
$row4H = $ws.Range("H4").Value2

$row7A = $ws.Range("A7").Value2
$row7B = $ws.Range("B7").Value2
$row7D = $ws.Range("D7").Value2
$row7E = $ws.Range("E7").Value2
$row7F = $ws.Range("F7").Value2
$row7G = $ws.Range("G7").Value2
$row7H = $ws.Range("H7").Value2

$ws.Range("A4").Value = $row7A
$ws.Range("B4").Value = $row7B
$ws.Range("D4").Value = $row7D
$ws.Range("E4").Value = $row7E
$ws.Range("F4").Value = $row7F
$ws.Range("G4").Value = $row7G
$ws.Range("H4").Value = $row7H

$ws.Range("A7").Value = $row4A
$ws.Range("B7").Value = $row4B
$ws.Range("D7").Value = $row4D
$ws.Range("E7").Value = $row4E
$ws.Range("F7").Value = $row4F
$ws.Range("G7").Value = $row4G
$ws.Range("H7").Value = $row4H

# --- Swap location data between row 5 and row 6 (columns A,Q,R) ---
$row5A = $ws.Range("A5").Value2
$row5Q = $ws.Range("Q5").Value2
$row5R = $ws.Range("R5").Value2

$row6A = $ws.Range("A6").Value2
$row6Q = $ws.Range("Q6").Value2
$row6R = $ws.Range("R6").Value2

$ws.Range("A5").Value = $row6A
$ws.Range("Q5").Value = $row6Q
$ws.Range("R5").Value = $row6R

$ws.Range("A6").Value = $row5A
$ws.Range("Q6").Value = $row5Q
$ws.Range("R6").Value = $row5R
